# Update TPM-derived NATMI ligand-receptor statistics on the active sheet.
# These cells recompute expression / specificity metrics after the TPM
# recalculation (per commit message: "update scripts wuth new tpm").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 12.367401
$ws.Range("H2").Value = 37.102203
$ws.Range("I2").Value = 0.8693283326101076
$ws.Range("J2").Value = 0.8693283326101076
$ws.Range("M2").Value = 1.492477333333333
$ws.Range("N2").Value = 4.477432
$ws.Range("O2").Value = 0.02769484181536182
$ws.Range("P2").Value = 0.02769484181536182
$ws.Range("Q2").Value = 18.458065664744
$ws.Range("R2").Value = 166.122590982696
$ws.Range("S2").Value = 0.02407591065724918
$ws.Range("T2").Value = 0.02407591065724918

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value = 12.367401
$ws.Range("H3").Value = 37.102203
$ws.Range("I3").Value = 0.8693283326101076
$ws.Range("J3").Value = 0.8693283326101076
$ws.Range("O3").Value = 0.6282762845978157
$ws.Range("P3").Value = 0.6282762845978156
$ws.Range("Q3").Value = 418.7337481117281
$ws.Range("R3").Value = 3768.603733005553
$ws.Range("S3").Value = 0.5461783749078926
$ws.Range("T3").Value = 0.5461783749078924

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value = 12.367401
$ws.Range("H4").Value = 37.102203
$ws.Range("I4").Value = 0.8693283326101076
$ws.Range("J4").Value = 0.8693283326101076
$ws.Range("N4").Value = 55.619234
$ws.Range("O4").Value = 0.3440288735868225
$ws.Range("P4").Value = 0.3440288735868225
$ws.Range("Q4").Value = 229.288456730278
$ws.Range("R4").Value = 2063.596110572502
$ws.Range("S4").Value = 0.2990740470449659
$ws.Range("T4").Value = 0.2990740470449659

# Row 5 (FAPs -> ECs)
$ws.Range("H5").Value = 4.303227
$ws.Range("I5").Value = 0.1008273593013545
$ws.Range("J5").Value = 0.1008273593013546
$ws.Range("M5").Value = 1.492477333333333
$ws.Range("N5").Value = 4.477432
$ws.Range("O5").Value = 0.02769484181536182
$ws.Range("P5").Value = 0.02769484181536182
$ws.Range("Q5").Value = 2.140822919229333
$ws.Range("R5").Value = 19.267406273064
$ws.Range("S5").Value = 0.002792397766511665
$ws.Range("T5").Value = 0.002792397766511665

# Row 6 (FAPs -> FAPs)
$ws.Range("H6").Value = 4.303227
$ws.Range("I6").Value = 0.1008273593013545
$ws.Range("J6").Value = 0.1008273593013546
$ws.Range("O6").Value = 0.6282762845978157
$ws.Range("P6").Value = 0.6282762845978156
$ws.Range("Q6").Value = 48.56602101728534
$ws.Range("R6").Value = 437.094189155568
$ws.Range("S6").Value = 0.06334743868766404
$ws.Range("T6").Value = 0.06334743868766404

# Row 7 (FAPs -> MuSCs)
$ws.Range("H7").Value = 4.303227
$ws.Range("I7").Value = 0.1008273593013545
$ws.Range("J7").Value = 0.1008273593013546
$ws.Range("N7").Value = 55.619234
$ws.Range("O7").Value = 0.3440288735868225
$ws.Range("P7").Value = 0.3440288735868225
$ws.Range("Q7").Value = 26.59357660756866
$ws.Range("R7").Value = 239.342189468118
$ws.Range("S7").Value = 0.03468752284717884
$ws.Range("T7").Value = 0.03468752284717885

# Row 8 (MuSCs -> ECs)
$ws.Range("I8").Value = 0.02984430808853782
$ws.Range("J8").Value = 0.02984430808853782
$ws.Range("M8").Value = 1.492477333333333
$ws.Range("N8").Value = 4.477432
$ws.Range("O8").Value = 0.02769484181536182
$ws.Range("P8").Value = 0.02769484181536182
$ws.Range("Q8").Value = 0.6336710512622222
$ws.Range("R8").Value = 5.70303946136
$ws.Range("S8").Value = 0.0008265333916009783
$ws.Range("T8").Value = 0.0008265333916009784

# Row 9 (MuSCs -> FAPs)
$ws.Range("I9").Value = 0.02984430808853782
$ws.Range("J9").Value = 0.02984430808853782
$ws.Range("O9").Value = 0.6282762845978157
$ws.Range("P9").Value = 0.6282762845978156
$ws.Range("S9").Value = 0.01875047100225908
$ws.Range("T9").Value = 0.01875047100225908

# Row 10 (MuSCs -> MuSCs)
$ws.Range("I10").Value = 0.02984430808853782
$ws.Range("J10").Value = 0.02984430808853782
$ws.Range("N10").Value = 55.619234
$ws.Range("O10").Value = 0.3440288735868225
$ws.Range("P10").Value = 0.3440288735868225
$ws.Range("S10").Value = 0.01026730369467776
$ws.Range("T10").Value = 0.01026730369467777
